# Update countries & provincias Spain
#
# The source "paises" data table (countries + COVID case counters) was
# refreshed with newer numbers for a handful of countries. Because the
# sheet is kept sorted descending by column B ("Casos totales"), those
# updated totals changed the relative ranking of a few neighbouring
# countries, which is why whole rows appear to "swap". The table below
# lists, for every row whose content changed, the final country name and
# the final B:H (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) values after the refresh+re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=31;  Name="Rumania";               Vals=@(7707, 491, 1357, 5958, 243, 20, 392)},
    @{Row=40;  Name="Indonesia";              Vals=@(5516, 380, 548,  4472,   0, 27, 496)},
    @{Row=75;  Name="Bosnia y Herzegovina";   Vals=@(1155,  45, 269,   843,   4,  2,  43)},
    @{Row=76;  Name="Lituania";               Vals=@(1128,  37, 178,   920,  14,  0,  30)},
    @{Row=79;  Name="Eslovaquia";             Vals=@(977,  114, 151,   818,   5,  2,   8)},
    @{Row=80;  Name="Republica de Macedonia"; Vals=@(974,    0,  98,   831,  15,  0,  45)},
    @{Row=102; Name="Malta";                  Vals=@(412,   13,  82,   327,   4,  0,   3)},
    @{Row=103; Name="Nigeria";                Vals=@(407,    0, 128,   267,   2,  0,  12)},
    @{Row=104; Name="Guinea";                 Vals=@(404,    0,  31,   372,   0,  0,   1)},
    @{Row=105; Name="Jordania";               Vals=@(401,    0, 250,   144,   5,  0,   7)},
    @{Row=185; Name="Suazilandia";            Vals=@(16,     1,   8,     8,   0,  0,   0)},
    @{Row=186; Name="Dominica";               Vals=@(16,     0,   8,     8,   0,  0,   0)},
    @{Row=196; Name="Islas Malvinas";         Vals=@(11,     0,   1,    10,   0,  0,   0)},
    @{Row=197; Name="Montserrat";             Vals=@(11,     0,   1,    10,   1,  0,   0)},
    @{Row=209; Name="Santo Tome y Principe";  Vals=@(4,      0,   0,     4,   0,  0,   0)},
    @{Row=210; Name="Sudan del Sur";          Vals=@(4,      0,   0,     4,   0,  0,   0)},
    @{Row=215; Name="Yemen";                  Vals=@(1,      0,   0,     1,   0,  0,   0)},
    @{Row=216; Name="San Pedro y Miquelon";   Vals=@(1,      0,   0,     1,   0,  0,   0)}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Name
    $col = 2
    foreach ($v in $u.Vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}
